$wb = $excel.ActiveWorkbook

# Remember the currently active sheet so we can restore it at the end
# (adding/copying sheets below will shift Excel's "active sheet" pointer).
$origActiveName = $wb.ActiveSheet.Name

# ------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new top data row for 2022-Q3,
#    pushing the existing quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 2.04

# The freshly inserted row copied formatting from the header row for
# B2:D2 - strip that back to the plain (unstyled) look used by the
# other data rows.
$summary.Range("B2:D2").ClearFormats()

# Column A on every data row carries the bold/bordered "index" style;
# restore it on the new row by copying the format from the row below.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Add the new "2022-Q3" detail sheet. It has the exact same shape
#    as the existing "2022-Q2" sheet, so clone that sheet (which also
#    places the clone right before it) and then update its figures.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2").Previous
$q3.Name = "2022-Q3"

# D2:G2 are stored as text in this workbook (e.g. "27.40", not 27.4),
# so use the text-prefix trick to keep them as text instead of being
# auto-coerced to numbers, then drop the leftover quote-prefix style.
$q3.Range("D2").Value = "'27.40"
$q3.Range("E2").Value = "'94.42"
$q3.Range("F2").Value = "'7.45"
$q3.Range("G2").Value = "'2.0413"
$q3.Range("D2:G2").ClearFormats()

$q3.Range("H2").Value = 8

# ------------------------------------------------------------------
# Restore whichever sheet was active before we started editing.
# ------------------------------------------------------------------
$wb.Worksheets.Item($origActiveName).Activate()
